$wb = $excel.ActiveWorkbook

# --- Rename the original sheet and add the new "R vs Stata" sheet right after it ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Testing summary"

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "R vs Stata"

# --- Populate the "R vs Stata" sheet -------------------------------------------------
$ws2.Range("A1").Value = "R vs Stata for the examples in the mimix help file"

$ws2.Range("C3").Value = "R results"
$ws2.Range("E3").Value = "Stata results"
$ws2.Range("G3").Value = "MC error"
$ws2.Range("I3").Value = "Z statistic, Stata - R"

$ws2.Range("C4").Value = "estimate"
$ws2.Range("D4").Value = "std.error"
$ws2.Range("E4").Value = "estimate"
$ws2.Range("F4").Value = "std.error"
$ws2.Range("G4").Value = "estimate"
$ws2.Range("H4").Value = "std.error"
$ws2.Range("I4").Value = "estimate"
$ws2.Range("J4").Value = "std.error"

$ws2.Range("B5").Value = "CC"
$ws2.Range("B6").Value = "MAR"
$ws2.Range("B7").Value = "LMCF"
$ws2.Range("B8").Value = "J2R"
$ws2.Range("B9").Value = "CIR"

# R results
$ws2.Range("C5").Value = 0.23854
$ws2.Range("D5").Value = 0.098540000000000003
$ws2.Range("C6").Value = 0.32844129999999999
$ws2.Range("D6").Value = 0.1033602
$ws2.Range("C7").Value = 0.29416930000000002
$ws2.Range("D7").Value = 0.098182350000000002
$ws2.Range("C8").Value = 0.22327810000000001
$ws2.Range("D8").Value = 0.10499171
$ws2.Range("C9").Value = 0.28442299999999998
$ws2.Range("D9").Value = 0.10153821

# Stata results
$ws2.Range("E5").Value = 0.23854075999999999
$ws2.Range("F5").Value = 0.098536529999999997
$ws2.Range("E6").Value = 0.32626279000000002
$ws2.Range("F6").Value = 0.10383219
$ws2.Range("E7").Value = 0.29572551000000002
$ws2.Range("F7").Value = 0.099239270000000004
$ws2.Range("E8").Value = 0.22448609
$ws2.Range("F8").Value = 0.10592255
$ws2.Range("E9").Value = 0.28327022000000002
$ws2.Range("F9").Value = 0.10248983

# MC error
$ws2.Range("G5").Value = 0
$ws2.Range("H5").Value = 0
$ws2.Range("G6").Value = 0.002297
$ws2.Range("H6").Value = 0.00094700000000000004
$ws2.Range("G7").Value = 0.0018779000000000001
$ws2.Range("H7").Value = 0.00064070000000000002
$ws2.Range("G8").Value = 0.0019273000000000001
$ws2.Range("H8").Value = 0.00067639999999999996
$ws2.Range("G9").Value = 0.0018341
$ws2.Range("H9").Value = 0.00068289999999999996

# Z statistic, Stata - R (shared formulas over I6:I9 and J6:J9)
$ws2.Range("I6:I9").Formula = "=(E6-C6)/(SQRT(2)*G6)"
$ws2.Range("J6:J9").Formula = "=(F6-D6)/(SQRT(2)*H6)"

# --- Number formats -------------------------------------------------------------------
$ws2.Range("C5:F9").NumberFormat = "0.0000"
$ws2.Range("G6:H9").NumberFormat = "0.0000"
$ws2.Range("I6:J9").NumberFormat = "0.00"

# --- Cosmetic sheet1 changes picked up by the newer Excel build ----------------------
$ws1.Rows.Item(2).RowHeight = 30
$ws1.Rows.Item(3).RowHeight = 105
$ws1.Rows.Item(4).RowHeight = 75
$ws1.Rows.Item(5).RowHeight = 30
$ws1.Rows.Item(6).RowHeight = 30
$ws1.Rows.Item(7).RowHeight = 60
$ws1.Rows.Item(8).RowHeight = 75
$ws1.Rows.Item(9).RowHeight = 60
$ws1.Rows.Item(10).RowHeight = 45
$ws1.Rows.Item(11).RowHeight = 105
$ws1.Rows.Item(12).RowHeight = 75
$ws1.Rows.Item(13).RowHeight = 45
$ws1.Rows.Item(14).RowHeight = 75
$ws1.Rows.Item(15).RowHeight = 45
$ws1.Rows.Item(16).RowHeight = 30
$ws1.Rows.Item(17).RowHeight = 45
$ws1.Rows.Item(18).RowHeight = 30
$ws1.Rows.Item(21).RowHeight = 30
$ws1.Rows.Item(22).RowHeight = 45
$ws1.Rows.Item(24).RowHeight = 30
$ws1.Rows.Item(25).RowHeight = 30
$ws1.Rows.Item(27).RowHeight = 45
$ws1.Rows.Item(29).RowHeight = 45
$ws1.Rows.Item(31).RowHeight = 30
$ws1.Rows.Item(34).RowHeight = 45
$ws1.Rows.Item(37).RowHeight = 45
$ws1.Rows.Item(38).RowHeight = 60
$ws1.Rows.Item(39).RowHeight = 60
$ws1.Rows.Item(40).RowHeight = 60
$ws1.Rows.Item(41).RowHeight = 75

# --- Select the new sheet as the active tab, matching the source workbook -------------
$ws2.Select()
$wb.Windows.Item(1).ScrollRow = 1
$wb.Windows.Item(1).ScrollColumn = 1
